$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: grupo/sexo/descripcion-ocupacion move from "iaest-dimension" to "iaest-measure"
$ws.Range("E2").Value = "iaest-measure:grupo"
$ws.Range("F2").Value = "iaest-measure:sexo"
$ws.Range("G2").Value = "iaest-measure:descripcion-ocupacion"

# Row 3: type changes from "dim" to "medida" for those same columns
$ws.Range("E3:G3").Value = "medida"

# Row 4: datatype changes from "skos:Concept" to "xsd:int" for those same columns
$ws.Range("E4:G4").Value = "xsd:int"

# Row 5: the per-column mapping file cells for grupo/sexo/descripcion-ocupacion are removed
# (measures, unlike curated dimensions, no longer carry a mapping-*.xlsx reference)
$ws.Range("E5:G5").Clear()
